# Applies the "Se crean las asociaciones de pokemon y location en hubspot" edit:
#  - adds an annotation cell F52 ("no lo esta guardando") to the Location block
#  - adds a new "Pokemon" association block (rows 66-76) mirroring the
#    existing Location/Move blocks, associating Pokemon properties with
#    HubSpot Contact properties

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Location block: flag that the "Number of areas" property isn't saving.
# ---------------------------------------------------------------------------
$ws.Range("F52").Borders.Item(10).LineStyle = 1
$ws.Range("F52").Value = "no lo está guardando"

# ---------------------------------------------------------------------------
# 2. New "Pokemon" association block (rows 66-76)
# ---------------------------------------------------------------------------

# Row 66: block title, styled like the "location" block title (row 49)
$ws.Range("B49:E49").Copy()
$ws.Range("B66:E66").PasteSpecial(-4122)
$ws.Range("B66:E66").Merge()
$ws.Range("B66").Value = "Pokemon"

# Row 67: column headers, styled like the "location"/"Move" headers (row 59)
# with the inner property-type columns styled like row 5
$ws.Range("B59:E59").Copy()
$ws.Range("B67:E67").PasteSpecial(-4122)
$ws.Range("C5:D5").Copy()
$ws.Range("C67:D67").PasteSpecial(-4122)
$ws.Rows.Item(67).RowHeight = 28.8
$ws.Range("B67").Value = "variable local"
$ws.Range("C67").Value = "Pokémon property"
$ws.Range("D67").Value = "HubSpot Contact"
$ws.Range("E67").Value = "variable Hubspot"

# Rows 68-75: body rows mapping each Pokemon property to its HubSpot
# counterpart, styled like the corresponding rows of the main Pokemon table
# (rows 6-13), but without borders on columns B/E (unlike that table).
$ws.Range("C6:D6").Copy()
$ws.Range("C68:D68").PasteSpecial(-4122)
$ws.Range("B68").Value = "id_pokemon"
$ws.Range("C68").Value = "id"
$ws.Range("D68").Value = "Pokédex id"
$ws.Range("E68").Value = "phone"

$ws.Range("C7:D7").Copy()
$ws.Range("C69:D69").PasteSpecial(-4122)
$ws.Range("B69").Value = "name"
$ws.Range("C69").Value = "Name"
$ws.Range("D69").Value = "name"
$ws.Range("E69").Value = "name"

$ws.Range("C8:D8").Copy()
$ws.Range("C70:D70").PasteSpecial(-4122)
$ws.Range("B70").Value = "hp"
$ws.Range("C70").Value = "hp"
$ws.Range("D70").Value = "HP"
$ws.Range("E70").Value = "hp"

$ws.Range("C9:D9").Copy()
$ws.Range("C71:D71").PasteSpecial(-4122)
$ws.Range("B71").Value = "attack"
$ws.Range("C71").Value = "attack"
$ws.Range("D71").Value = "Attack"
$ws.Range("E71").Value = "attack"

$ws.Range("C10:D10").Copy()
$ws.Range("C72:D72").PasteSpecial(-4122)
$ws.Range("B72").Value = "defense"
$ws.Range("C72").Value = "defense"
$ws.Range("D72").Value = "Defense"
$ws.Range("E72").Value = "defense"

$ws.Range("C11:D11").Copy()
$ws.Range("C73:D73").PasteSpecial(-4122)
$ws.Rows.Item(73).RowHeight = 28.8
$ws.Range("B73").Value = "special-defense"
$ws.Range("C73").Value = "special-defense"
$ws.Range("D73").Value = "Special Defense"
$ws.Range("E73").Value = "special_defense"

$ws.Range("C12:D12").Copy()
$ws.Range("C74:D74").PasteSpecial(-4122)
$ws.Rows.Item(74).RowHeight = 28.8
$ws.Range("B74").Value = "special-attack"
$ws.Range("C74").Value = "special-attack"
$ws.Range("D74").Value = "Special Attack"
$ws.Range("E74").Value = "special_attack"

$ws.Range("C13:D13").Copy()
$ws.Range("C75:D75").PasteSpecial(-4122)
$ws.Range("B75").Value = "speed"
$ws.Range("C75").Value = "speed"
$ws.Range("D75").Value = "Speed"
$ws.Range("E75").Value = "speed"

# Row 76: last body row (types), plus a "tabla" marker in column A and a
# "<----- multi check" annotation in column F.
$ws.Range("C14:D14").Copy()
$ws.Range("C76:D76").PasteSpecial(-4122)
$ws.Range("A76").Value = "tabla"
$ws.Range("B76").Value = "type"
$ws.Range("C76").Value = "types"
$ws.Range("D76").Value = "Types"
$ws.Range("E76").Value = "type"
$ws.Range("F76").Value = "<----- multi check"

# ---------------------------------------------------------------------------
# 3. Restore selection to match the authored state.
# ---------------------------------------------------------------------------
$ws.Range("H82").Select()
